$wb = $excel.ActiveWorkbook

# --- Sheet references (before any structural changes) ------------------
$total = $wb.Worksheets.Item(1)   # 总计
$q3    = $wb.Worksheets.Item(2)   # 2022-Q3 (keeps its original figures,
                                   # simply becomes "2022-Q3" again once the
                                   # new Q4 sheet is inserted ahead of it)

# --- Insert the new "2022-Q4" sheet right after 总计 -------------------
# Copying the existing 2022-Q3 sheet keeps all of its formatting (styles,
# column layout, page margins, etc.) and places the new sheet directly
# before the sheet we copied from, i.e. right after 总计.
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Update the Q4 sheet's figures (fund code/name stay the same). These are
# stored as text in the source workbook, so force text ('-prefix) instead
# of letting numeric-looking strings be auto-converted to numbers.
$q4.Range("D2").Value = "'0.39"
$q4.Range("E2").Value = "'94.72"
$q4.Range("F2").Value = "'2.22"
$q4.Range("G2").Value = "'0.0087"
$q4.Range("H2").Value = 10

# --- Update the summary ("总计") sheet ----------------------------------
# Row 2 now documents 2022-Q4, row 3 documents 2022-Q3 (shifted down from
# what used to be row 2), and a new row 4 documents 2022-Q2 - duplicated
# from row 3 (same format) with updated values.
$total.Range("B2").Value = "2022-Q4"
$total.Range("B3").Value = "2022-Q3"

$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)  # xlPasteFormats
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.01

# --- Restore the original active tab (2022-Q2) --------------------------
$q2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q2.Activate()
